# Update countries & provincias Spain
# This script applies the data refresh captured in the commit:
#  - Swap the display order of "Hungria"/"Serbia" (rows 76-77) and
#    "Georgia"/"Maldivas" (rows 105-106), because the underlying ranking by
#    total cases changed.
#  - Refresh the numeric COVID counters for the affected rows.
#  - Update the "Datos actualizados..." timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 28 (Ucrania) ---
$ws.Range("B28").Value = 250538
$ws.Range("C28").Value = 5804
$ws.Range("D28").Value = 110650
$ws.Range("E28").Value = 135109
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 89
$ws.Range("H28").Value = 4779

# --- Row 63 (Armenia) ---
$ws.Range("B63").Value = 55087
$ws.Range("C63").Value = 614
$ws.Range("D63").Value = 45528
$ws.Range("E63").Value = 8549
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 1010

# --- Row 75 (Afganistan) ---
$ws.Range("B75").Value = 39693
$ws.Range("C75").Value = 77
$ws.Range("D75").Value = 33058
$ws.Range("E75").Value = 5163
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 1472

# --- Row 76: was Serbia, now Hungria (new data) ---
$ws.Range("A76").Value = "Hungria"
$ws.Range("B76").Value = 35222
$ws.Range("C76").Value = 1176
$ws.Range("D76").Value = 9202
$ws.Range("E76").Value = 25107
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 15
$ws.Range("H76").Value = 913

# --- Row 77: was Hungria, now Serbia (Serbia's previous data) ---
$ws.Range("A77").Value = "Serbia"
$ws.Range("B77").Value = 34344
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 31536
$ws.Range("E77").Value = 2048
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 760

# --- Row 79 (El Salvador) ---
$ws.Range("E79").Value = 4191
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 881

# --- Row 105: was Maldivas, now Georgia (new data) ---
$ws.Range("A105").Value = "Georgia"
$ws.Range("B105").Value = 10752
$ws.Range("C105").Value = 527
$ws.Range("D105").Value = 5866
$ws.Range("E105").Value = 4814
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 6
$ws.Range("H105").Value = 72

# --- Row 106: was Georgia, now Maldivas (Maldivas' previous data) ---
$ws.Range("A106").Value = "Maldivas"
$ws.Range("B106").Value = 10742
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 9589
$ws.Range("E106").Value = 1119
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 34

# --- Row 176 (Taiwan) ---
$ws.Range("B176").Value = 527
$ws.Range("C176").Value = 3
$ws.Range("D176").Value = 488

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 09:17"
